# "Loan RBI, Variable Instalments"
#
# The Repayment Schedule sheet gains a new (blank) column between the
# existing "In Advance" (M) and "Late" (N) columns, which pushes the old
# N/O/P columns ("Late", "Outstanding"/"Heading", "Over Due") one slot to
# the right (-> O/P/Q). The sheet also becomes the active tab/selection
# (previously it was "Transactions").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before the old "Late" column (N), shifting
# the old N/O/P columns to O/P/Q.
$ws.Columns("N").Insert()

# Excel gave the freshly inserted column the same width as its neighbour
# to the left ("In Advance", column M).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment Schedule" the active sheet/tab again (it becomes
# activeTab in the workbook, and Transactions loses tabSelected), with
# the cursor left on R8.
$ws.Activate()
$ws.Range("R8").Select() | Out-Null
